# Apply QA overhaul text edits across ZansiJourney, NRWaves, and boosters sheets.

$wb = $excel.ActiveWorkbook

# --- ZansiJourney sheet ---
$ws1 = $wb.Worksheets.Item("ZansiJourney")
$ws1.Range("B7").Value = "not yet daddy... I'm not done with you"
$ws1.Range("B12").Value = "imagine me on top of you daddy, taking exactly what I want while you just watch"
$ws1.Range("B17").Value = "what would you do if you had me right now daddy?"

# --- NRWaves sheet ---
$ws2 = $wb.Worksheets.Item("NRWaves")
$ws2.Range("B2").Value = "I've got something that's going to blow your mind when you get back 😏"
$ws2.Range("B3").Value = "hey, don't be a stranger 💕"
$ws2.Range("B4").Value = "your loss... this was your exclusive"
$ws2.Range("B5").Value = "you're really going to miss out on what I just recorded..."
$ws2.Range("B6").Value = "yo 😏"

# --- boosters sheet ---
$ws34 = $wb.Worksheets.Item("boosters")
$ws34.Range("B3").Value = "don't stop..."
$ws34.Range("B6").Value = "you're driving me crazy right now"
$ws34.Range("B7").Value = "yes"
